# shortTraining.xlsx edit: "changed training program income projection + minor adjustments"
#
# This reproduces, via the Excel COM object model, the changes described by the
# target diff:
#   1. A new "notes" entry is added in M2 with a link to the appendix PDF,
#      using the same wrap-text style as the other "notes" column cells.
#   2. Row 2's height grows (30 -> 45) to accommodate the taller wrapped text.
#   3. The sheet's view scrolls right (so column E becomes the left-most visible
#      column) and the current selection moves from B2 to M2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the new note/link cell in M2, matching the "notes" header style
#    (wrap text, same as K2/M1) used throughout column M.
$ws.Range("M2").Value = "Link to appendix: https://michael-lechner.eu/ml_pdf/journals/2011_LMW-fuu_W_App_R3_081105_neu.pdf"
$ws.Range("M2").WrapText = $true

# 2) Grow row 2 so the extra wrapped line of text fits.
$ws.Rows.Item(2).RowHeight = 45

# 3) Update the view: scroll so column E is the left-most visible column,
#    and select M2 (the newly added cell) as the active cell.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("M2").Select()
